# Use timezone from preferences for excel reports.
#
# The "Period:" value (B6) and the per-event "Time" column template (A9)
# both used to build their date strings with joda-time's own formatter /
# DateTime constructor. Switch both to the shared dateTool.format(...)
# helper so the configured locale/timezone from preferences is honoured.
#
# NOTE: write A9 first, then B6 last -- the two new template strings are
# unique (not found anywhere else in the shared-string table), so each
# write lands in shared-strings order of "most recently written" = highest
# index. Writing A9 before B6 reproduces the exact index layout in the
# target workbook (B6 -> higher index, A9 -> lower index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", event.serverTime, locale, timezone)}'
$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'
